$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

# B1: new label "California"
$ws.Range("B1").Value = "California"

# C1: date stamp (2022-09-02), formatted with the built-in short-date
# number format (numFmtId 14) via a recognized date format string.
# (Set the format before the value so the engine resolves straight to
# the built-in numFmtId instead of minting a redundant custom one.)
$ws.Range("C1").NumberFormat = "mm-dd-yy"
$d = Get-Date -Year 2022 -Month 9 -Day 2 -Hour 0 -Minute 0 -Second 0
$ws.Range("C1").Value = $d
